$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.011618333333334
$ws.Range("H2").Value = 3.034855
$ws.Range("I2").Value = 0.5235149663433657
$ws.Range("J2").Value = 0.5235149663433657
$ws.Range("M2").Value = 1.522526333333333
$ws.Range("N2").Value = 4.567579
$ws.Range("O2").Value = 0.2115373313282365
$ws.Range("P2").Value = 0.2115373313282365
$ws.Range("Q2").Value = 1.540215551782778
$ws.Range("R2").Value = 13.861939966045
$ws.Range("S2").Value = 0.1107429588906671
$ws.Range("T2").Value = 0.1107429588906671

# Row 3
$ws.Range("G3").Value = 1.011618333333334
$ws.Range("H3").Value = 3.034855
$ws.Range("I3").Value = 0.5235149663433657
$ws.Range("J3").Value = 0.5235149663433657
$ws.Range("O3").Value = 0.4376697219060474
$ws.Range("P3").Value = 0.4376697219060474
$ws.Range("Q3").Value = 3.186698574627223
$ws.Range("R3").Value = 28.68028717164501
$ws.Range("S3").Value = 0.2291266497331546
$ws.Range("T3").Value = 0.2291266497331546

# Row 4
$ws.Range("G4").Value = 1.011618333333334
$ws.Range("H4").Value = 3.034855
$ws.Range("I4").Value = 0.5235149663433657
$ws.Range("J4").Value = 0.5235149663433657
$ws.Range("M4").Value = 2.524809666666667
$ws.Range("N4").Value = 7.574429
$ws.Range("O4").Value = 0.3507929467657161
$ws.Range("P4").Value = 0.3507929467657162
$ws.Range("Q4").Value = 2.554143746977223
$ws.Range("R4").Value = 22.987293722795
$ws.Range("S4").Value = 0.1836453577195439
$ws.Range("T4").Value = 0.183645357719544

# Row 5
$ws.Range("I5").Value = 0.2899264353016711
$ws.Range("J5").Value = 0.2899264353016712
$ws.Range("M5").Value = 1.522526333333333
$ws.Range("N5").Value = 4.567579
$ws.Range("O5").Value = 0.2115373313282365
$ws.Range("P5").Value = 0.2115373313282365
$ws.Range("Q5").Value = 0.8529826905305555
$ws.Range("R5").Value = 7.676844214775
$ws.Range("S5").Value = 0.06133026440522412
$ws.Range("T5").Value = 0.06133026440522413

# Row 6
$ws.Range("I6").Value = 0.2899264353016711
$ws.Range("J6").Value = 0.2899264353016712
$ws.Range("O6").Value = 0.4376697219060474
$ws.Range("P6").Value = 0.4376697219060474
$ws.Range("S6").Value = 0.1268920223116941
$ws.Range("T6").Value = 0.1268920223116941

# Row 7
$ws.Range("I7").Value = 0.2899264353016711
$ws.Range("J7").Value = 0.2899264353016712
$ws.Range("M7").Value = 2.524809666666667
$ws.Range("N7").Value = 7.574429
$ws.Range("O7").Value = 0.3507929467657161
$ws.Range("P7").Value = 0.3507929467657162
$ws.Range("Q7").Value = 1.414503575669444
$ws.Range("R7").Value = 12.730532181025
$ws.Range("S7").Value = 0.1017041485847529
$ws.Range("T7").Value = 0.101704148584753

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.360498
$ws.Range("H8").Value = 1.081494
$ws.Range("I8").Value = 0.1865585983549632
$ws.Range("J8").Value = 0.1865585983549632
$ws.Range("M8").Value = 1.522526333333333
$ws.Range("N8").Value = 4.567579
$ws.Range("O8").Value = 0.2115373313282365
$ws.Range("P8").Value = 0.2115373313282365
$ws.Range("Q8").Value = 0.548867698114
$ws.Range("R8").Value = 4.939809283026
$ws.Range("S8").Value = 0.03946410803234524
$ws.Range("T8").Value = 0.03946410803234524

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.360498
$ws.Range("H9").Value = 1.081494
$ws.Range("I9").Value = 0.1865585983549632
$ws.Range("J9").Value = 0.1865585983549632
$ws.Range("O9").Value = 0.4376697219060474
$ws.Range("P9").Value = 0.4376697219060474
$ws.Range("Q9").Value = 1.135604629634
$ws.Range("R9").Value = 10.220441666706
$ws.Range("S9").Value = 0.08165104986119874
$ws.Range("T9").Value = 0.08165104986119874

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.360498
$ws.Range("H10").Value = 1.081494
$ws.Range("I10").Value = 0.1865585983549632
$ws.Range("J10").Value = 0.1865585983549632
$ws.Range("M10").Value = 2.524809666666667
$ws.Range("N10").Value = 7.574429
$ws.Range("O10").Value = 0.3507929467657161
$ws.Range("P10").Value = 0.3507929467657162
$ws.Range("Q10").Value = 0.910188835214
$ws.Range("R10").Value = 8.191699516926001
$ws.Range("S10").Value = 0.06544344046141921
$ws.Range("T10").Value = 0.06544344046141923
